{"js": "// Update the two-digit division worksheet: replace each \"A\u00f7B=\" prompt\n// with a new \"A\u00f7B=\" prompt, cell by cell (row, col), preserving all\n// paragraph/run formatting (fonts, size, alignment) already on the runs.\n//\n// The worksheet is a single 5-column table; every 4th row (0, 4, 8, 12, 16)\n// holds the five division problems for that line, the rows in between are\n// blank spacer rows. We address each cell by (rowIndex, colIndex) rather\n// than searching by text, because several prompts repeat verbatim\n// (e.g. \"38\u00f72=\" and \"59\u00f74=\" each occur twice) and position-based lookup\n// removes any ambiguity about which occurrence to touch.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> [ [oldText, newText], ... ] for the 5 columns in that row\nconst rowEdits = {\n  0: [\n    [\"37\u00f78=\", \"12\u00f76=\"],\n    [\"38\u00f72=\", \"60\u00f73=\"],\n    [\"79\u00f72=\", \"51\u00f72=\"],\n    [\"29\u00f72=\", \"99\u00f77=\"],\n    [\"71\u00f73=\", \"61\u00f76=\"],\n  ],\n  4: [\n    [\"83\u00f78=\", \"26\u00f73=\"],\n    [\"60\u00f78=\", \"33\u00f77=\"],\n    [\"77\u00f74=\", \"20\u00f75=\"],\n    [\"39\u00f75=\", \"29\u00f77=\"],\n    [\"42\u00f76=\", \"21\u00f72=\"],\n  ],\n  8: [\n    [\"45\u00f79=\", \"59\u00f73=\"],\n    [\"59\u00f74=\", \"82\u00f78=\"],\n    [\"36\u00f79=\", \"21\u00f79=\"],\n    [\"38\u00f72=\", \"33\u00f79=\"],\n    [\"37\u00f73=\", \"70\u00f79=\"],\n  ],\n  12: [\n    [\"47\u00f79=\", \"76\u00f75=\"],\n    [\"13\u00f75=\", \"33\u00f73=\"],\n    [\"38\u00f77=\", \"54\u00f79=\"],\n    [\"59\u00f74=\", \"22\u00f76=\"],\n    [\"19\u00f79=\", \"40\u00f72=\"],\n  ],\n  16: [\n    [\"35\u00f78=\", \"23\u00f75=\"],\n    [\"21\u00f74=\", \"73\u00f74=\"],\n    [\"39\u00f74=\", \"48\u00f76=\"],\n    [\"10\u00f78=\", \"72\u00f72=\"],\n    [\"62\u00f72=\", \"36\u00f73=\"],\n  ],\n};\n\n// First load all the cell paragraph texts so we can sanity-check before\n// mutating (defensive: avoid clobbering a cell whose text unexpectedly\n// doesn't match what the diff expects).\nconst cellRefs = [];\nfor (const rowIndexStr of Object.keys(rowEdits)) {\n  const rowIndex = Number(rowIndexStr);\n  const edits = rowEdits[rowIndex];\n  for (let colIndex = 0; colIndex < edits.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    cellRefs.push({ rowIndex, colIndex, para, edits: edits[colIndex] });\n  }\n}\nawait context.sync();\n\nfor (const ref of cellRefs) {\n  const [oldText, newText] = ref.edits;\n  if (ref.para.text !== oldText) {\n    throw new Error(\n      `Unexpected text in cell (${ref.rowIndex},${ref.colIndex}): ` +\n        `expected \"${oldText}\" but found \"${ref.para.text}\"`\n    );\n  }\n  ref.para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit division worksheet: replace each \"A\u00f7B=\" prompt\n# with a new \"A\u00f7B=\" prompt, cell by cell (row, column), preserving all\n# paragraph/run formatting already present in the table.\n#\n# The worksheet is a single 5-column table; every 4th row (1, 5, 9, 13, 17\n# in Word's 1-based Table.Cell indexing) holds the five division problems\n# for that line, the rows in between are blank spacer rows. Cells are\n# addressed by (row, column) rather than located via Find/Replace, because\n# several prompts repeat verbatim (e.g. \"38\u00f72=\" and \"59\u00f74=\" each occur\n# twice) and position-based lookup removes any ambiguity about which\n# occurrence to touch.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# rowIndex (1-based) -> list of [oldText, newText] for columns 1..5\n$rowEdits = @{\n    1  = @(\n        @(\"37\u00f78=\", \"12\u00f76=\"),\n        @(\"38\u00f72=\", \"60\u00f73=\"),\n        @(\"79\u00f72=\", \"51\u00f72=\"),\n        @(\"29\u00f72=\", \"99\u00f77=\"),\n        @(\"71\u00f73=\", \"61\u00f76=\")\n    )\n    5  = @(\n        @(\"83\u00f78=\", \"26\u00f73=\"),\n        @(\"60\u00f78=\", \"33\u00f77=\"),\n        @(\"77\u00f74=\", \"20\u00f75=\"),\n        @(\"39\u00f75=\", \"29\u00f77=\"),\n        @(\"42\u00f76=\", \"21\u00f72=\")\n    )\n    9  = @(\n        @(\"45\u00f79=\", \"59\u00f73=\"),\n        @(\"59\u00f74=\", \"82\u00f78=\"),\n        @(\"36\u00f79=\", \"21\u00f79=\"),\n        @(\"38\u00f72=\", \"33\u00f79=\"),\n        @(\"37\u00f73=\", \"70\u00f79=\")\n    )\n    13 = @(\n        @(\"47\u00f79=\", \"76\u00f75=\"),\n        @(\"13\u00f75=\", \"33\u00f73=\"),\n        @(\"38\u00f77=\", \"54\u00f79=\"),\n        @(\"59\u00f74=\", \"22\u00f76=\"),\n        @(\"19\u00f79=\", \"40\u00f72=\")\n    )\n    17 = @(\n        @(\"35\u00f78=\", \"23\u00f75=\"),\n        @(\"21\u00f74=\", \"73\u00f74=\"),\n        @(\"39\u00f74=\", \"48\u00f76=\"),\n        @(\"10\u00f78=\", \"72\u00f72=\"),\n        @(\"62\u00f72=\", \"36\u00f73=\")\n    )\n}\n\nforeach ($rowIndex in $rowEdits.Keys) {\n    $edits = $rowEdits[$rowIndex]\n    for ($i = 0; $i -lt $edits.Count; $i++) {\n        $colIndex = $i + 1\n        $oldText = $edits[$i][0]\n        $newText = $edits[$i][1]\n\n        $cell = $tbl.Cell($rowIndex, $colIndex)\n        $rng = $cell.Range\n\n        # Cell range text includes the trailing end-of-cell marks\n        # (carriage return + cell marker); strip them for comparison.\n        $current = $rng.Text.TrimEnd([char]13, [char]7)\n\n        if ($current -ne $oldText) {\n            throw \"Unexpected text in cell ($rowIndex,$colIndex): expected '$oldText' but found '$current'\"\n        }\n\n        $rng.Text = $newText\n    }\n}\n"}
